$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "new"

$newSheet.Range("A1").Value = "Outdoor Model"
$newSheet.Range("B1").Value = "Outdoor Quantity"
$newSheet.Range("C1").Value = "Outdoor Serial(s)"
$newSheet.Range("D1").Value = "Indoor Model"
$newSheet.Range("E1").Value = "Indoor Quantity"
$newSheet.Range("F1").Value = "Indoor Serial(s)"
